$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")
$ws.Activate()

# Add new player row (Week 15 simulations) for N.Bellore with all zero stats
$ws.Range("A8").Value = "N.Bellore"
$ws.Range("B8:J8").Value = 0

# Leave selection where Excel would land after entering the row of data
$ws.Range("G9").Select()
